# Insert two new data rows ("Betarraga" price records) above the former
# row 237, pushing the existing rows 237:291 down to 239:293. Then
# populate the two new rows with their full record data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 237-291 down by two rows.
$ws.Rows("237:238").Insert()

# --- New row 237 ---
$ws.Range("A237").Value = 10
$ws.Range("B237").Value = "Vega Modelo de Temuco"
$ws.Range("C237").Value = "La Araucanía"
$ws.Range("D237").Value = 44511
$ws.Range("E237").Value = 9
$ws.Range("F237").Value = 100114014
$ws.Range("G237").Value = "Betarraga"
$ws.Range("H237").Value = "Sin especificar"
$ws.Range("I237").Value = "Primera"
$ws.Range("J237").Value = 40
$ws.Range("K237").Value = 10000
$ws.Range("L237").Value = 10000
$ws.Range("M237").Value = 10000
$ws.Range("N237").Value = "`$/docena de paquetes"
$ws.Range("O237").Value = "Provincia de Cautín"
$ws.Range("P237").Value = 833
$ws.Range("Q237").Value = 12
$ws.Range("R237").Value = "Hortaliza"

# --- New row 238 ---
$ws.Range("A238").Value = 10
$ws.Range("B238").Value = "Vega Modelo de Temuco"
$ws.Range("C238").Value = "La Araucanía"
$ws.Range("D238").Value = 44511
$ws.Range("E238").Value = 9
$ws.Range("F238").Value = 100114014
$ws.Range("G238").Value = "Betarraga"
$ws.Range("H238").Value = "Sin especificar"
$ws.Range("I238").Value = "Primera"
$ws.Range("J238").Value = 200
$ws.Range("K238").Value = 700
$ws.Range("L238").Value = 700
$ws.Range("M238").Value = 700
$ws.Range("N238").Value = "`$/paquete 5 unidades"
$ws.Range("O238").Value = "Región Metropolitana"
$ws.Range("P238").Value = 140
$ws.Range("Q238").Value = 5
$ws.Range("R238").Value = "Hortaliza"
